$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.411.91"
$ws.Range("E2").Value = "  +2.63%  "

$ws.Range("D3").Value = "1.842.15"
$ws.Range("E3").Value = "  +1.88%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").Value = "'229.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.10%  "

$ws.Range("D6").Value = "'0.609"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.49%  "

$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("D8").Value = "'43.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.96%  "

$ws.Range("E9").Value = "  +6.91%  "

$ws.Range("D10").Value = "'0.0696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.54%  "

$ws.Range("E11").Value = "  +3.84%  "

$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").Value = "1.850.14"
$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("D14").Value = "'11.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("E15").Value = "  +7.15%  "

$ws.Range("E16").Value = "  +6.86%  "

$ws.Range("D17").Value = "35.368.29"
$ws.Range("E17").Value = "  +2.65%  "

$ws.Range("D18").Value = "'70.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.63%  "

$ws.Range("E19").Value = "  +3.47%  "

$ws.Range("D20").Value = "'244.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("D21").Value = "'12.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.32%  "

$ws.Range("D22").Value = "'4.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.56%  "

$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").Value = "'168.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("D26").Value = "'7.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.54%  "

$ws.Range("D27").Value = "'17.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("E29").Value = "  +13.70%  "

$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("D31").Value = "3.346.58"
$ws.Range("E31").Value = "  +37.74%  "

$ws.Range("E32").Value = "  +6.31%  "

$ws.Range("D33").Value = "'4.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.76%  "

$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("E35").Value = "  +2.30%  "

$ws.Range("D36").Value = "'95.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +15.09%  "

$ws.Range("E37").Value = "  +7.65%  "

$ws.Range("D38").Value = "1.344.07"
$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("E39").Value = "  +2.81%  "

$ws.Range("E40").Value = "  +3.46%  "

$ws.Range("D41").Value = "'2.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.94%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'15.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.38%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.30%  "

$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0519"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.33%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'6.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.33%  "

$ws.Range("D49").Value = "2.009.38"
$ws.Range("E49").Value = "  +2.03%  "

$ws.Range("E50").Value = "  +0.32%  "

$ws.Range("D51").Value = "'102.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
